# "Generate Report for Archive"
#
# The c9df683d-f4c3-46aa-ad8f-78ec726b97c6.md file has moved from
# "Ready for handoff" to "In Translation", and its handback information
# (Latest Target File / Latest Handback File / Latest Handback DateTime)
# has been recorded for the archive report on both the zh-cn and de-de
# language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: update the status column for the c9df683d row
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B5").Value = "In Translation"
$overview.Range("C5").Value = "In Translation"

# ---------------------------------------------------------------
# Helper language-sheet update: zh-cn and de-de both get the same
# shape of change on row 5 (the c9df683d file).
# ---------------------------------------------------------------
function Update-LangSheet($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status -> In Translation
    $ws.Range("B5").Value = "In Translation"

    # Find the existing hyperlink addresses for A5 (source md) and C5
    # (latest handoff xlf) so the new E5/F5 hyperlinks point at the same
    # targets (this mirrors the reviewer re-running the same source file
    # through the report generator for the handback columns).
    $addrA5 = ""
    $addrC5 = ""
    $dispA5 = ""
    $dispC5 = ""
    foreach ($h in $ws.Hyperlinks) {
        $r = $h.Range.Address()
        if ($r -eq '$A$5') {
            $addrA5 = $h.Address
            $dispA5 = $h.TextToDisplay
        }
        if ($r -eq '$C$5') {
            $addrC5 = $h.Address
            $dispC5 = $h.TextToDisplay
        }
    }

    # Latest Target File (E5) and Latest Handback File (F5)
    $ws.Range("E5").Value = $dispA5
    $ws.Range("F5").Value = $dispC5
    $ws.Range("E5").Style = "HyperLink"
    $ws.Range("F5").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("E5"), $addrA5, "", "", $dispA5) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F5"), $addrC5, "", "", $dispC5) | Out-Null

    # Latest Handback DateTime (G5) - now recorded
    $ws.Range("G5").Value = $handbackDateTime
}

Update-LangSheet "zh-cn" "2016-01-25 06:02:29"
Update-LangSheet "de-de" "2016-01-25 06:02:46"
